# New eSM test data: update the flow-cytometry channel names used in the
# process_fcs(...) transformation formula (FSC/SSC/FL1 -> FSC-H/SSC-H/FL1-H)
# and refresh the per-sheet active-cell selections left over from editing.

$wb = $excel.ActiveWorkbook

# --- Samples sheet: just move the selection ---
$wsSamples = $wb.Worksheets.Item("Samples")
$wsSamples.Range("C7").Select()

# --- Groups sheet: just move the selection ---
$wsGroups = $wb.Worksheets.Item("Groups")
$wsGroups.Range("A2").Select()

# --- Transformations sheet: update the process_fcs formula text and move the selection ---
$wsTransformations = $wb.Worksheets.Item("Transformations")
$wsTransformations.Range("B2").Value = 'process_fcs("plate_01",["FSC-H","SSC-H"],["FL1-H"])'
$wsTransformations.Range("E13").Select()

# --- Views sheet: just move the selection ---
$wsViews = $wb.Worksheets.Item("Views")
$wsViews.Range("A3").Select()

# Leave the Samples sheet active/focused (matches tabSelected="1" on sheet1).
$wsSamples.Activate()
